$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows (2-9) ---
$ws.Cells.Item(2,3).Value = "ts1-2.jpg"
$ws.Cells.Item(3,3).Value = "ts1-3.jpg"
$ws.Cells.Item(4,3).Value = "ts1-4.jpg"

$ws.Cells.Item(5,2).Value = 2
$ws.Cells.Item(5,3).Value = "ts4-2.jpg"

$ws.Cells.Item(6,2).Value = 2
$ws.Cells.Item(6,3).Value = "ts4-3.jpg"

$ws.Cells.Item(7,2).Value = 2
$ws.Cells.Item(7,3).Value = "ts4-4.jpg"

$ws.Cells.Item(8,2).Value = 3
$ws.Cells.Item(8,3).Value = "sl1-2.jpg"

$ws.Cells.Item(9,2).Value = 3
$ws.Cells.Item(9,3).Value = "sl1-3.jpg"

# --- Append new rows (10-30) ---
$data = @(
    @(9,  3,  "sl1-4.jpg"),
    @(10, 4,  "hd1-2.jpg"),
    @(11, 4,  "hd1-3.jpg"),
    @(12, 4,  "hd1-4.jpg"),
    @(13, 5,  "hd2-2.jpg"),
    @(14, 5,  "hd2-3.jpg"),
    @(15, 5,  "hd2-4.jpg"),
    @(16, 6,  "tas1-2.jpg"),
    @(17, 6,  "tas1-3.jpg"),
    @(18, 6,  "tas1-4.jpg"),
    @(19, 7,  "tas3-1.jpg"),
    @(20, 7,  "tas3-2.jpg"),
    @(21, 7,  "tas3-3.jpg"),
    @(22, 8,  "sc1-2.jpg"),
    @(23, 8,  "sc1-3.jpg"),
    @(24, 9,  "sr1-2.jpg"),
    @(25, 9,  "sr1-3.jpg"),
    @(26, 9,  "sr1-4.jpg"),
    @(27, 10, "js1-2.jpg"),
    @(28, 10, "js1-3.jpg"),
    @(29, 10, "js1-4.jpg")
)

$row = 10
foreach ($entry in $data) {
    $ws.Cells.Item($row,1).Value = $entry[0]
    $ws.Cells.Item($row,2).Value = $entry[1]
    $ws.Cells.Item($row,3).Value = $entry[2]
    $row = $row + 1
}

# --- Selection / view state ---
$ws.Range("F29").Select()
